# Append the 11/30/2025 row of profit data (row 6) to the sheet,
# mirroring the existing rows 2-5.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces the date to be stored as plain text (matching
# the existing A2:A5 date cells) instead of being auto-converted to a date
# serial number. Resetting the style back to Normal afterwards avoids
# picking up a stray text/quote-prefix number format on the cell.
$ws.Range("A6").Value = "'11/30/2025"
$ws.Range("A6").Style = "Normal"

$ws.Range("B6").Value = 14780.29
$ws.Range("C6").Value = 0.1573518266434222
$ws.Range("D6").Value = 0.8426481733565778
$ws.Range("E6").Value = -46.38
$ws.Range("F6").Value = -11.31
$ws.Range("G6").Value = -18012.17
$ws.Range("H6").Value = -59.12
$ws.Range("I6").Value = -477.14
$ws.Range("J6").Value = -17.02
